# Edit applied to cooccurrence/post_objects.xlsx
#
# Summary of the content change (per the supplied diff):
#  - Sheet1, row 11 ("chair" / "table") is removed entirely, shrinking the
#    used range from A1:B72 to A1:B71 and shifting every later row up by one.
#  - The small "target -> landmark" lookup table in column B (rows 1-7) is
#    refreshed: the kitchen-appliance words (cabinet, cupboard, fridge,
#    grill, oven, range, console, island, refrigerator) are replaced with
#    the furniture/room words (counter, floor, table, shelf, wall, stove)
#    that remain relevant after the "chair" row was dropped.
#  - The leftover B8:B14 cells (old grill/oven/range/console/island/
#    refrigerator/stove tail) are cleared since the new lookup table only
#    spans rows 1-7.
#  - The active selection moves from B16 to B6.
#
# (Sheet2 keeps the exact same 8 text values - only the shared-string
# indices they point at shift around, which falls out automatically once
# the workbook is re-saved.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "chair" / "table" row - shifts A12:B72 up to A11:B71.
$ws.Rows.Item(11).Delete()

# Refresh the lookup table in column B (header B1 "landmark" is unchanged).
$ws.Range("B2").Value = "counter"
$ws.Range("B3").Value = "floor"
$ws.Range("B4").Value = "table"
$ws.Range("B5").Value = "shelf"
$ws.Range("B6").Value = "wall"
$ws.Range("B7").Value = "stove"

# Drop the now-unused tail of the old lookup table (full Clear so the
# cells disappear from the saved XML instead of lingering as empty-but-
# styled placeholders).
$ws.Range("B8:B14").Clear()

# Match the saved selection state (activeCell B6).
[void]$ws.Range("B6").Select()
